$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the title text for row 4 (C34) to include the {movie_id} placeholder
$ws.Range("B4").Value = "Read movie by id with existing {movie_id}"

# Remove row 10 (the duplicate TestRail-exported "C1" case) and shift rows up
$ws.Rows.Item(10).Delete()
